$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": a new day ("22-dec") is inserted as a column right
#     before the existing "01-oct." block, shifting everything from ES
#     onward one column to the right (ES:FW -> ET:FX).
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Columns.Item(149).Insert(-4161)
$ws1.Range("ES1").Value = "22-dec"
$ws1.Range("ES2:ES25").Value = "-"

# --- Sheet "Gaz": append the two newest daily prices.
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A177").NumberFormat = "@"
$ws2.Range("A177").Value = "2025-12-20"
$ws2.Range("A177").Style = "Normal"
$ws2.Range("B177").Value = 26.9
$ws2.Range("A178").NumberFormat = "@"
$ws2.Range("A178").Value = "2025-12-21"
$ws2.Range("A178").Style = "Normal"
$ws2.Range("B178").Value = 26.9

# --- Sheet "CO2": append the two newest daily prices.
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A178").NumberFormat = "@"
$ws3.Range("A178").Value = "2025-12-20"
$ws3.Range("A178").Style = "Normal"
$ws3.Range("B178").Value = 84.54000000000001
$ws3.Range("A179").NumberFormat = "@"
$ws3.Range("A179").Value = "2025-12-21"
$ws3.Range("A179").Style = "Normal"
$ws3.Range("B179").Value = 84.54000000000001
